$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 currently holds phone "09876543" as text. The edit splits this into
# two rows: row 24 becomes the numeric phone 9876543 (points reset to 0),
# and a new row 25 is appended that keeps the original leading-zero text
# value "09876543" (also with 0 points).
#
# Copy row 24 down into the newly-needed row 25 first (this preserves the
# original text-typed phone value and the numeric 0 points), then overwrite
# row 24's phone with the plain numeric value.
$ws.Range("A24:C24").Copy()
$ws.Range("A25").PasteSpecial()

$ws.Range("A24").Value = 9876543
